$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates per the crypto price/volume refresh.
# Price cells (column D) that would otherwise be auto-coerced to numbers
# by Excel (losing trailing zeros / the text type) are written with a
# leading apostrophe to force text, then restyled back to Normal so no
# stray "quote prefix" cell style is left behind.

$ws.Range("D2").Value = "68.428.96"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.649.23"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'596.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'158.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "2.648.18"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'0.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'28.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "3.132.68"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'0.0000188"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "68.291.91"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "2.630.43"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "'11.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "'364.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'7.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'4.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'74.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'9.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "2.799.94"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "'568.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'8.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").Value = "'1.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'160.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "'1.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'5.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'2.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "0.0₆0321"
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("D46").Value = "'158.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "'3.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "'21.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "'0.0780"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").Value = "'0.575"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.09%  "
